$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column A (date-like text "11-04-2023") as Text so Excel does not
# auto-convert it to a date serial number; reset style afterwards to avoid
# leaving a residual number-format style on the cells.
$ws.Range("A2:A10").NumberFormat = "@"

$ws.Range("A2").Value = "11-04-2023"
$ws.Range("B2").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C2").Value = "10/04/2023  08:44:11"
$ws.Range("D2").Value = 12161017
$ws.Range("E2").Value = 8597595
$ws.Range("F2").Value = "BRAYAN MIGUEL JAFFRA PEREIRA"
$ws.Range("G2").Value = "10 dias úteis"
$ws.Range("H2").Value = "NO"
$ws.Range("I2").Value = "No Assistencial"
$ws.Range("J2").Value = "Responder  Detalhes"

$ws.Range("A3").Value = "11-04-2023"
$ws.Range("B3").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C3").Value = "10/04/2023  11:56:29"
$ws.Range("D3").Value = 12161650
$ws.Range("E3").Value = 8598342
$ws.Range("F3").Value = "JANAINA MARIA DA SILVA"
$ws.Range("G3").Value = "10 dias úteis"
$ws.Range("H3").Value = "NO"
$ws.Range("I3").Value = "Assistencial"
$ws.Range("J3").Value = "Responder  Detalhes"

$ws.Range("A4").Value = "11-04-2023"
$ws.Range("B4").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C4").Value = "10/04/2023  13:25:58"
$ws.Range("D4").Value = 12162010
$ws.Range("E4").Value = 8598761
$ws.Range("F4").Value = "CLEO DALSIOR VOM DOELINGER"
$ws.Range("G4").Value = "10 dias úteis"
$ws.Range("H4").Value = "NO"
$ws.Range("I4").Value = "Assistencial"
$ws.Range("J4").Value = "Responder  Detalhes"

$ws.Range("A5").Value = "11-04-2023"
$ws.Range("B5").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C5").Value = "10/04/2023  14:54:09"
$ws.Range("D5").Value = 12162415
$ws.Range("E5").Value = 8599264
$ws.Range("F5").Value = "WELLINGTON FERREIRA DE JESUS"
$ws.Range("G5").Value = "10 dias úteis"
$ws.Range("H5").Value = "NO"
$ws.Range("I5").Value = "Assistencial"
$ws.Range("J5").Value = "Responder  Detalhes"

$ws.Range("A6").Value = "11-04-2023"
$ws.Range("B6").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C6").Value = "10/04/2023  15:00:01"
$ws.Range("D6").Value = 12162439
$ws.Range("E6").Value = 8599230
$ws.Range("F6").Value = "BENICIO MARTINS FERNANDES"
$ws.Range("G6").Value = "10 dias úteis"
$ws.Range("H6").Value = "NO"
$ws.Range("I6").Value = "Assistencial"
$ws.Range("J6").Value = "Responder  Detalhes"

$ws.Range("A7").Value = "11-04-2023"
$ws.Range("B7").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C7").Value = "11/04/2023  08:23:00"
$ws.Range("D7").Value = 12163370
$ws.Range("E7").Value = 8600412
$ws.Range("F7").Value = "PAOLA MONIQUE DA SILVA TEIXEIRA"
$ws.Range("G7").Value = "10 dias úteis"
$ws.Range("H7").Value = "NO"
$ws.Range("I7").Value = "Assistencial"
$ws.Range("J7").Value = "Responder  Detalhes"

$ws.Range("A8").Value = "11-04-2023"
$ws.Range("B8").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C8").Value = "11/04/2023  09:34:28"
$ws.Range("D8").Value = 12163535
$ws.Range("E8").Value = 8600621
$ws.Range("F8").Value = "IDEILDO LUCENA MOURA DA SILVA JUNIOR"
$ws.Range("G8").Value = "10 dias úteis"
$ws.Range("H8").Value = "NO"
$ws.Range("I8").Value = "Assistencial"
$ws.Range("J8").Value = "Responder  Detalhes"

$ws.Range("A9").Value = "11-04-2023"
$ws.Range("B9").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C9").Value = "11/04/2023  13:50:53"
$ws.Range("D9").Value = 12164502
$ws.Range("E9").Value = 8601769
$ws.Range("F9").Value = "MARIA DA CONCEICAO DA SILVA"
$ws.Range("G9").Value = "10 dias úteis"
$ws.Range("H9").Value = "NO"
$ws.Range("I9").Value = "Assistencial"
$ws.Range("J9").Value = "Responder  Detalhes"

$ws.Range("A10").Value = "11-04-2023"
$ws.Range("B10").Value = "422380 - YOU ASSISTÊNCIA MÉDICA LTDA"
$ws.Range("C10").Value = "11/04/2023  15:08:53"
$ws.Range("D10").Value = 12164852
$ws.Range("E10").Value = 8602234
$ws.Range("F10").Value = "VALDIR FERNANDES DE ARAUJO"
$ws.Range("G10").Value = "10 dias úteis"
$ws.Range("H10").Value = "NO"
$ws.Range("I10").Value = "Assistencial"
$ws.Range("J10").Value = "Responder  Detalhes"

# Reset style back to Normal so only the value/type differs from the original
# (matches upstream workbook which keeps these cells unstyled).
$ws.Range("A2:A10").Style = "Normal"
